# [FEATURE] agregar los mismos campos al importar que en el formulario
#
# The "Medio" column (old col G, with values TELÉFONO/WHATSAPP/CORREO/
# FORMULARIO PÁGINA WEB/REDES) is removed, and three new columns -
# Dirección, Responsable, Ciudad - are inserted right before the old
# "Fuente" column (old col F), so the customer-contact columns that used
# to follow (Nombres/Teléfonos/Emails/Cargos Contactos Cliente) shift two
# slots to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "Medio" column entirely (old column G).
$ws.Columns("G").Delete()

# Insert three fresh blank columns for Dirección / Responsable / Ciudad
# right before the "Fuente" column (old column F, now still F after the
# delete above).
$ws.Range("F1:H1").EntireColumn.Insert()

# New header row labels for the inserted columns.
$ws.Range("F1").Value = "Dirección"
$ws.Range("G1").Value = "Responsable"
$ws.Range("H1").Value = "Ciudad"

# Sample data for the first customer row, matching the webform fields.
$ws.Range("F2").Value = "Avenida Siempre Viva 123"
$ws.Range("G2").Value = "Jhonatan"
$ws.Range("H2").Value = "Envigado"

# Column widths for the new columns (matches manual resize in Excel).
$ws.Columns("F").ColumnWidth = 40.7109375
$ws.Range("G1:H1").EntireColumn.ColumnWidth = 25.5703125

# Stray underline formatting left on the empty F13 cell (artifact of the
# author's manual edit session).
$ws.Range("F13").Font.Underline = $true

# Restore the active selection/cell seen in the saved file.
$ws.Range("H22").Select()
